$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.323.87'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.332.11'
$ws.Range('E3').Value = '  -0.30%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '188.95'
$ws.Range('E5').Value = '  +3.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '561.33'
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.324.36'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.184'
$ws.Range('E10').Value = '  -0.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.591'
$ws.Range('E11').Value = '  +0.71%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.90'
$ws.Range('E12').Value = '  +1.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000273'
$ws.Range('E13').Value = '  +2.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.70'
$ws.Range('E14').Value = '  +1.38%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.864.00'
$ws.Range('E15').Value = '  -0.26%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '606.74'
$ws.Range('E16').Value = '  +0.48%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '66.402.58'
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.10'
$ws.Range('E18').Value = '  -0.12%  '
$ws.Range('E19').Value = '  +1.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.339.05'
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.17'
$ws.Range('E21').Value = '  -2.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.917'
$ws.Range('E22').Value = '  +1.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '18.61'
$ws.Range('E23').Value = '  +10.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.16'
$ws.Range('E24').Value = '  +1.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '100.59'
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('E26').Value = '  -0.60%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.00'
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('E28').Value = '  +4.41%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.72'
$ws.Range('E29').Value = '  +4.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.70'
$ws.Range('E30').Value = '  -0.16%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.57'
$ws.Range('E31').Value = '  -0.24%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.82'
$ws.Range('E32').Value = '  +8.96%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.89'
$ws.Range('E33').Value = '  +3.58%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '582.41'
$ws.Range('E34').Value = '  +8.79%  '
$ws.Range('E35').Value = '  +1.54%  '
$ws.Range('E36').Value = '  +1.65%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.725.30'
$ws.Range('E37').Value = '  -2.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '57.26'
$ws.Range('E38').Value = '  -1.09%  '
$ws.Range('E39').Value = '  +0.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0733'
$ws.Range('E40').Value = '  +2.59%  '
$ws.Range('E41').Value = '  +5.58%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '34.08'
$ws.Range('E42').Value = '  +6.91%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.47'
$ws.Range('E43').Value = '  +12.30%  '
$ws.Range('E44').Value = '  -5.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.71'
$ws.Range('E45').Value = '  +1.45%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.346'
$ws.Range('E46').Value = '  +1.58%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0427'
$ws.Range('E47').Value = '  +3.17%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.35'
$ws.Range('E48').Value = '  +3.47%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.130'
$ws.Range('E49').Value = '  +1.01%  '
$ws.Range('E50').Value = '  +0.08%  '
